$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the "Area" calculation columns
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Row 2: first segment area uses 0 as the lower depth bound
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
# Running total of the area column (segments 2 through 11)
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Row 3: area of the second segment
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15 share the same relative formula pattern
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Match the saved selection shown in the target file
[void]$ws.Range("F2").Select()
